$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slovenia Prva Liga")

# Row 78
$ws.Range("B78").Value = 5499423
$ws.Range("F78").Value = "Olimpija Ljubljana"
$ws.Range("G78").Value = "NK Celje"
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 2
$ws.Range("J78").Value = "A"
$ws.Range("K78").Value = 2.5
$ws.Range("L78").Value = 3.3
$ws.Range("M78").Value = 2.5
$ws.Range("N78").Value = 2.55
$ws.Range("O78").Value = 3.25
$ws.Range("P78").Value = 2.45
$ws.Range("Q78").Value = 0
$ws.Range("T78").Value = 2.5
$ws.Range("U78").Value = 1.85
$ws.Range("V78").Value = 1.95
$ws.Range("X78").Value = -1
$ws.Range("Y78").Value = 1.45
$ws.Range("AC78").Value = 0.95

# Row 79
$ws.Range("B79").Value = 5498504
$ws.Range("F79").Value = "NK Maribor"
$ws.Range("G79").Value = "NK Bravo"
$ws.Range("H79").Value = 1
$ws.Range("I79").Value = 1
$ws.Range("J79").Value = "D"
$ws.Range("K79").Value = 1.571
$ws.Range("L79").Value = 3.8
$ws.Range("M79").Value = 4.75
$ws.Range("N79").Value = 1.533
$ws.Range("O79").Value = 4
$ws.Range("P79").Value = 4.75
$ws.Range("Q79").Value = -1
$ws.Range("T79").Value = 2.75
$ws.Range("U79").Value = 1.875
$ws.Range("V79").Value = 1.925
$ws.Range("X79").Value = 3
$ws.Range("Y79").Value = -1
$ws.Range("AC79").Value = 0.925

# Row 80
$ws.Range("B80").Value = 5495053
$ws.Range("F80").Value = "NK Radomlje"
$ws.Range("G80").Value = "NK Domzale"
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = "H"
$ws.Range("K80").Value = 2.55
$ws.Range("L80").Value = 3.1
$ws.Range("M80").Value = 2.55
$ws.Range("N80").Value = 3.75
$ws.Range("P80").Value = 1.833
$ws.Range("Q80").Value = 0.5
$ws.Range("R80").Value = 1.925
$ws.Range("S80").Value = 1.875
$ws.Range("U80").Value = 1.975
$ws.Range("V80").Value = 1.825
$ws.Range("W80").Value = 2.75
$ws.Range("Y80").Value = -1
$ws.Range("Z80").Value = 0.925
$ws.Range("AA80").Value = -1
$ws.Range("AB80").Value = -1
$ws.Range("AC80").Value = 0.825

# Row 81
$ws.Range("B81").Value = 5498503
$ws.Range("F81").Value = "FC Koper"
$ws.Range("G81").Value = "NS Mura"
$ws.Range("I81").Value = 2
$ws.Range("J81").Value = "A"
$ws.Range("K81").Value = 2.05
$ws.Range("L81").Value = 3.3
$ws.Range("M81").Value = 3.25
$ws.Range("N81").Value = 2
$ws.Range("P81").Value = 3.25
$ws.Range("Q81").Value = -0.5
$ws.Range("R81").Value = 2
$ws.Range("S81").Value = 1.8
$ws.Range("U81").Value = 1.825
$ws.Range("V81").Value = 1.975
$ws.Range("W81").Value = -1
$ws.Range("Y81").Value = 2.25
$ws.Range("Z81").Value = -1
$ws.Range("AA81").Value = 0.8
$ws.Range("AB81").Value = 0.825
$ws.Range("AC81").Value = -1

# Row 186
$ws.Range("H186").Value = 4
$ws.Range("I186").Value = 1
$ws.Range("J186").Value = "H"
$ws.Range("N186").Value = 1.222
$ws.Range("O186").Value = 6.5
$ws.Range("P186").Value = 9
$ws.Range("Q186").Value = -2
$ws.Range("R186").Value = 2.025
$ws.Range("S186").Value = 1.775
$ws.Range("U186").Value = 1.825
$ws.Range("V186").Value = 1.975
$ws.Range("W186").Value = 0.222
$ws.Range("X186").Value = -1
$ws.Range("Y186").Value = -1
$ws.Range("Z186").Value = 1.025
$ws.Range("AA186").Value = -1
$ws.Range("AB186").Value = 0.825
$ws.Range("AC186").Value = -1

# Row 187
$ws.Range("H187").Value = 2
$ws.Range("I187").Value = 1
$ws.Range("J187").Value = "H"
$ws.Range("N187").Value = 1.571
$ws.Range("P187").Value = 5
$ws.Range("Q187").Value = -0.75
$ws.Range("R187").Value = 1.8
$ws.Range("S187").Value = 2
$ws.Range("U187").Value = 1.85
$ws.Range("V187").Value = 1.95
$ws.Range("W187").Value = 0.571
$ws.Range("X187").Value = -1
$ws.Range("Y187").Value = -1
$ws.Range("Z187").Value = 0.4
$ws.Range("AA187").Value = -0.5
$ws.Range("AB187").Value = 0.8500000000000001
$ws.Range("AC187").Value = -1

# Row 188
$ws.Range("N188").Value = 5.75
$ws.Range("P188").Value = 1.45
$ws.Range("Q188").Value = 1.25
$ws.Range("R188").Value = 1.775
$ws.Range("S188").Value = 2.025

# Row 189
$ws.Range("R189").Value = 1.925
$ws.Range("S189").Value = 1.875
$ws.Range("T189").Value = 2.75
$ws.Range("U189").Value = 1.8
$ws.Range("V189").Value = 2

# Row 190
$ws.Range("N190").Value = 6.5
$ws.Range("O190").Value = 4.75
$ws.Range("P190").Value = 1.363
$ws.Range("Q190").Value = 1.25
$ws.Range("R190").Value = 1.975
$ws.Range("S190").Value = 1.825
$ws.Range("U190").Value = 1.85
$ws.Range("V190").Value = 1.95

# Row 191
$ws.Range("N191").Value = 2.2
$ws.Range("R191").Value = 1.975
$ws.Range("S191").Value = 1.825
$ws.Range("U191").Value = 1.95
$ws.Range("V191").Value = 1.85
